# Updates cryptocurrency Price (D) and Volume(1h) (E) columns to refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.006.23"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "'3.733.99"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'601.77"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'167.94"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").Value = "'3.732.36"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +4.50%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'4.357.72"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "'3.738.64"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'68.946.74"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +6.85%  "
$ws.Range("D21").Value = "'497.77"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").Value = "'10.23"
$ws.Range("E22").Value = "  +13.59%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "'85.30"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").Value = "'2.43"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "'31.70"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").Value = "'3.886.27"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "'3.665.36"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'436.04"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").Value = "'48.99"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("D48").Value = "'40.52"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "'141.90"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "'2.744.09"
$ws.Range("E51").Value = "  -1.67%  "
